$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds values that look numeric (e.g. "1.040", "0.4622")
# but must stay as literal text strings, exactly as scraped (trailing
# zeros / multi-dot "thousands" notation would otherwise be silently
# coerced into floating point numbers by COM's auto-type-detection).
# Force text format only on the individual cells whose price actually
# changes, right before writing the new value, so untouched cells keep
# their original (default) style.

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.075.39"
$ws.Range("E2").Value = "  -2.47%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.901.98"
$ws.Range("E3").Value = "  -3.20%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -1.20%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.96"
$ws.Range("E5").Value = "  +1.21%  "

# Row 6 - USDC (price unchanged, volume changed)
$ws.Range("E6").Value = "  -1.19%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4622"
$ws.Range("E7").Value = "  -2.88%  "

# Row 8 - Cardano (price unchanged, volume changed)
$ws.Range("E8").Value = "  -1.58%  "

# Row 9 - OKB
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.94"
$ws.Range("E9").Value = "  -3.25%  "

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08311"
$ws.Range("E10").Value = "  -2.19%  "

# Row 11 - Polygon
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.040"
$ws.Range("E11").Value = "  -1.66%  "

# Row 12 - Solana
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.81"
$ws.Range("E12").Value = "  -3.07%  "

# Row 13 - WrappedEther
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.903.65"
$ws.Range("E13").Value = "  -4.66%  "

# Row 14 - Chainlink
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.350"
$ws.Range("E14").Value = "  -3.66%  "

# Row 15 - Polkadot
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.014"
$ws.Range("E15").Value = "  -3.62%  "

# Row 16 - BinanceUSD
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  -1.10%  "

# Row 17 - Litecoin
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.42"
$ws.Range("E17").Value = "  -0.77%  "

# Row 18 - ShibaInu (price unchanged, volume changed)
$ws.Range("E18").Value = "  -0.15%  "

# Row 19 - TRON
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06583"
$ws.Range("E19").Value = "  -0.45%  "

# Row 20 - Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.74"
$ws.Range("E20").Value = "  -4.76%  "

# Row 21 - Dai
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  -1.03%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.678"
$ws.Range("E22").Value = "  -1.78%  "

# Row 23 - WrappedBTC
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.059.33"
$ws.Range("E23").Value = "  -2.57%  "

# Row 24 - Cosmos
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.14"
$ws.Range("E24").Value = "  -3.45%  "

# Row 25 - Toncoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.314"
$ws.Range("E25").Value = "  +0.75%  "

# Row 26 - WrappedliquidstakedEther2.0
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.127.61"
$ws.Range("E26").Value = "  -4.55%  "

# Row 27 - Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.44"
$ws.Range("E27").Value = "  +0.51%  "

# Row 28 - EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.01"
$ws.Range("E28").Value = "  -0.84%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.121"
$ws.Range("E29").Value = "  -1.74%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.701"
$ws.Range("E30").Value = "  -4.47%  "

# Row 31 - BitcoinCash
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.15"
$ws.Range("E31").Value = "  +0.79%  "

# Row 32 - Stellar
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09588"
$ws.Range("E32").Value = "  -0.35%  "

# Row 33 - ImmutableX
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9671"
$ws.Range("E33").Value = "  -3.54%  "

# Row 34 - ARBITRUM
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.463"
$ws.Range("E34").Value = "  +0.17%  "

# Row 35 - HuobiToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.626"
$ws.Range("E35").Value = "  -1.89%  "

# Row 36 - Filecoin
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.490"
$ws.Range("E36").Value = "  -3.05%  "

# Row 37 - TrustWalletToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.261"
$ws.Range("E37").Value = "  -0.83%  "

# Row 38 - VeChain (price unchanged, volume changed)
$ws.Range("E38").Value = "  -2.32%  "

# Row 39 - FraxShare
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.673"
$ws.Range("E39").Value = "  -0.87%  "

# Row 40 - Hedera
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06122"
$ws.Range("E40").Value = "  -1.91%  "

# Row 41 - TheSandbox
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6111"
$ws.Range("E41").Value = "  -1.96%  "

# Rows 42/43 swap places: Frax moves up to 42, Aptos moves down to 43,
# each carrying its own refreshed price/volume.
$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("E42").Value = "  -1.12%  "

$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.81"
$ws.Range("E43").Value = "  -2.85%  "

# Row 44 - Algorand
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1899"
$ws.Range("E44").Value = "  -0.92%  "

# Row 45 - WEMIXTOKEN
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.308"
$ws.Range("E45").Value = "  -1.73%  "

# Row 46 - Decentraland
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5840"
$ws.Range("E46").Value = "  -1.77%  "

# Row 47 - EnergySwap
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.74"
$ws.Range("E47").Value = "  -0.99%  "

# Row 48 - NEARProtocol
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.002"
$ws.Range("E48").Value = "  -4.00%  "

# Row 49 - PancakeSwap
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.429"
$ws.Range("E49").Value = "  +0.14%  "

# Row 50 - Cronos
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06902"
$ws.Range("E50").Value = "  +0.87%  "

# Row 51 - Quant
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.78"
$ws.Range("E51").Value = "  +0.23%  "
